$d = $word.ActiveDocument

# Locate the anchor paragraph: "Lead comprehensive research initiatives..."
$anchorText = "Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions"

$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
        break
    }
    $i = $i + 1
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph not found"
}

$newParas = @(
    "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters",
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

# Paragraphs collection is 1-based; the anchor paragraph's Item index is anchorIndex + 1
$insertAfterItem = $anchorIndex + 1

foreach ($text in $newParas) {
    $p = $d.Paragraphs.Item($insertAfterItem)
    $p.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($insertAfterItem + 1)
    $newP.Range.Text = $text
    $insertAfterItem = $insertAfterItem + 1
}

Write-Output "Inserted $($newParas.Length) paragraphs after index $anchorIndex"
